$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds the same date value (serial 46061) for every
# data row (rows 2-412). The update bumps that date forward by one day
# (serial 46062) for every row.
for ($r = 2; $r -le 412; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 46061) {
        $cell.Value2 = 46062
    }
}
